$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Determine the last used row based on column A
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

# Swap the contents of column C (codeforiati:group-name) and column D
# (codeforiati:group-code), including the header row, for every used row.
for ($r = 1; $r -le $lastRow; $r++) {
    $cCell = $ws.Cells.Item($r, 3)
    $dCell = $ws.Cells.Item($r, 4)

    $cValue = $cCell.Value2
    $dValue = $dCell.Value2

    $cCell.Value = $dValue
    $dCell.Value = $cValue
}
